$wb = $excel.ActiveWorkbook
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "trend"

$headers = @(
    "ANCOM-BC2 (No Filter).30, 0.1 (N = 100)",
    "ANCOM-BC2 (SS Filter).30, 0.1 (N = 100)",
    "ANCOM-BC2 (No Filter).30, 0.2 (N = 100)",
    "ANCOM-BC2 (SS Filter).30, 0.2 (N = 100)",
    "ANCOM-BC2 (No Filter).30, 0.5 (N = 100)",
    "ANCOM-BC2 (SS Filter).30, 0.5 (N = 100)",
    "ANCOM-BC2 (No Filter).30, 0.9 (N = 100)",
    "ANCOM-BC2 (SS Filter).30, 0.9 (N = 100)",
    "ANCOM-BC2 (No Filter).60, 0.1 (N = 100)",
    "ANCOM-BC2 (SS Filter).60, 0.1 (N = 100)",
    "ANCOM-BC2 (No Filter).60, 0.2 (N = 100)",
    "ANCOM-BC2 (SS Filter).60, 0.2 (N = 100)",
    "ANCOM-BC2 (No Filter).60, 0.5 (N = 100)",
    "ANCOM-BC2 (SS Filter).60, 0.5 (N = 100)",
    "ANCOM-BC2 (No Filter).60, 0.9 (N = 100)",
    "ANCOM-BC2 (SS Filter).60, 0.9 (N = 100)",
    "ANCOM-BC2 (No Filter).90, 0.1 (N = 100)",
    "ANCOM-BC2 (SS Filter).90, 0.1 (N = 100)",
    "ANCOM-BC2 (No Filter).90, 0.2 (N = 100)",
    "ANCOM-BC2 (SS Filter).90, 0.2 (N = 100)",
    "ANCOM-BC2 (No Filter).90, 0.5 (N = 100)",
    "ANCOM-BC2 (SS Filter).90, 0.5 (N = 100)",
    "ANCOM-BC2 (No Filter).90, 0.9 (N = 100)",
    "ANCOM-BC2 (SS Filter).90, 0.9 (N = 100)",
    "ANCOM-BC2 (No Filter).150, 0.1 (N = 100)",
    "ANCOM-BC2 (SS Filter).150, 0.1 (N = 100)",
    "ANCOM-BC2 (No Filter).150, 0.2 (N = 100)",
    "ANCOM-BC2 (SS Filter).150, 0.2 (N = 100)",
    "ANCOM-BC2 (No Filter).150, 0.5 (N = 100)",
    "ANCOM-BC2 (SS Filter).150, 0.5 (N = 100)",
    "ANCOM-BC2 (No Filter).150, 0.9 (N = 100)",
    "ANCOM-BC2 (SS Filter).150, 0.9 (N = 100)",
    "ANCOM-BC2 (No Filter).300, 0.1 (N = 100)",
    "ANCOM-BC2 (SS Filter).300, 0.1 (N = 100)",
    "ANCOM-BC2 (No Filter).300, 0.2 (N = 100)",
    "ANCOM-BC2 (SS Filter).300, 0.2 (N = 100)",
    "ANCOM-BC2 (No Filter).300, 0.5 (N = 100)",
    "ANCOM-BC2 (SS Filter).300, 0.5 (N = 100)",
    "ANCOM-BC2 (No Filter).300, 0.9 (N = 100)",
    "ANCOM-BC2 (SS Filter).300, 0.9 (N = 100)"
)
for ($c = 0; $c -lt $headers.Count; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

$row2 = @(1, 0.99, 1, 0.98, 0.96, 0.92, 0.49, 0.46, 1, 1, 1, 0.99, 0.98, 0.94, 0.52, 0.5, 1, 1, 1, 0.99, 0.97, 0.93, 0.54, 0.51, 1, 1, 1, 1, 0.98, 0.94, 0.54, 0.52, 1, 1, 1, 1, 0.99, 0.97, 0.55, 0.52)
for ($c = 0; $c -lt $row2.Count; $c++) {
    $ws.Cells.Item(2, $c + 1).Value = $row2[$c]
}

$row3 = @(0, 0.03, 0, 0.02, 0.09, 0.11, 0.05, 0.05, 0, 0.01, 0, 0.02, 0.04, 0.07, 0.02, 0.02, 0, 0, 0, 0.01, 0.06, 0.07, 0.01, 0.01, 0, 0, 0, 0.01, 0.03, 0.05, 0.01, 0.01, 0, 0, 0, 0, 0.03, 0.04, 0, 0.01)
for ($c = 0; $c -lt $row3.Count; $c++) {
    $ws.Cells.Item(3, $c + 1).Value = $row3[$c]
}

$row4 = @(0.03, 0.02, 0.02, 0.01, 0.09, 0.07, 0, 0, 0.02, 0.01, 0.01, 0, 0.12, 0.11, 0, 0, 0.01, 0.01, 0.01, 0, 0.08, 0.07, 0, 0, 0.01, 0, 0, 0, 0.06, 0.06, 0, 0, 0.01, 0, 0, 0, 0.03, 0.03, 0, 0)
for ($c = 0; $c -lt $row4.Count; $c++) {
    $ws.Cells.Item(4, $c + 1).Value = $row4[$c]
}

$row5 = @(0.1, 0.09, 0.06, 0.06, 0.12, 0.1, 0, 0, 0.06, 0.05, 0.03, 0.03, 0.16, 0.15, 0, 0, 0.08, 0.07, 0.04, 0.04, 0.16, 0.15, 0, 0, 0.04, 0.03, 0.01, 0.01, 0.14, 0.14, 0, 0, 0.02, 0.02, 0.01, 0.01, 0.12, 0.12, 0, 0)
for ($c = 0; $c -lt $row5.Count; $c++) {
    $ws.Cells.Item(5, $c + 1).Value = $row5[$c]
}
